# Apply crypto price/volume updates to match the target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Text)
    $cell = $ws.Range($CellRef)
    # Force text storage (no numeric/date auto-coercion), then restore the
    # default "Normal" style so no visible formatting change is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "58.164.16"
$ws.Range("E2").Value = "  -3.76%  "
$ws.Range("D3").Value = "3.139.20"
$ws.Range("E3").Value = "  -5.21%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue "D5" "523.58"
$ws.Range("E5").Value = "  -6.35%  "
Set-TextValue "D6" "134.87"
$ws.Range("E6").Value = "  -5.33%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.137.93"
$ws.Range("E8").Value = "  -5.22%  "
$ws.Range("E9").Value = "  -5.92%  "
$ws.Range("E10").Value = "  -7.47%  "
$ws.Range("E11").Value = "  -8.59%  "
$ws.Range("E12").Value = "  -6.61%  "
$ws.Range("D13").Value = "3.675.86"
$ws.Range("E13").Value = "  -5.30%  "
$ws.Range("E14").Value = "  -1.49%  "
Set-TextValue "D15" "25.51"
$ws.Range("E15").Value = "  -5.13%  "
$ws.Range("D16").Value = "3.138.64"
$ws.Range("E16").Value = "  -5.36%  "
$ws.Range("D17").Value = "58.110.02"
$ws.Range("E17").Value = "  -3.86%  "
Set-TextValue "D18" "0.0000153"
$ws.Range("E18").Value = "  -7.54%  "
Set-TextValue "D19" "5.84"
$ws.Range("E19").Value = "  -4.99%  "
Set-TextValue "D20" "13.05"
$ws.Range("E20").Value = "  -8.96%  "
Set-TextValue "D21" "7.97"
$ws.Range("E21").Value = "  -8.29%  "
Set-TextValue "D22" "345.23"
$ws.Range("E22").Value = "  -7.81%  "
$ws.Range("E23").Value = "  +0.07%  "
Set-TextValue "D24" "68.65"
$ws.Range("E24").Value = "  -8.30%  "
Set-TextValue "D25" "0.507"
$ws.Range("E25").Value = "  -6.02%  "
$ws.Range("D26").Value = "3.263.35"
$ws.Range("E26").Value = "  -5.35%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D27" "0.168"
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0956"
$ws.Range("E28").Value = "  -6.88%  "
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D30" "6.78"
$ws.Range("E30").Value = "  -5.73%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D31" "0.998"
$ws.Range("E31").Value = "  -0.11%  "
Set-TextValue "D32" "1.86"
$ws.Range("E32").Value = "  -9.15%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D33" "6.83"
$ws.Range("E33").Value = "  -10.12%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D34" "21.51"
$ws.Range("E34").Value = "  -5.06%  "
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D36" "4.80"
$ws.Range("E36").Value = "  -6.92%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D37" "157.26"
$ws.Range("E37").Value = "  -5.52%  "
Set-TextValue "D38" "6.23"
$ws.Range("E38").Value = "  -7.28%  "
$ws.Range("E39").Value = "  -11.02%  "
Set-TextValue "D40" "0.0693"
$ws.Range("E40").Value = "  -5.19%  "
$ws.Range("D41").Value = "3.169.90"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D42" "24.32"
$ws.Range("E42").Value = "  -9.94%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "40.46"
$ws.Range("E43").Value = "  -3.69%  "
Set-TextValue "D44" "0.695"
$ws.Range("E44").Value = "  -7.59%  "
Set-TextValue "D45" "1.08"
$ws.Range("E45").Value = "  -3.09%  "
Set-TextValue "D46" "3.91"
$ws.Range("E46").Value = "  -6.17%  "
Set-TextValue "D47" "1.00"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  -9.05%  "
$ws.Range("D49").Value = "2.265.50"
$ws.Range("E49").Value = "  -4.59%  "
Set-TextValue "D50" "6.19"
$ws.Range("E50").Value = "  -3.79%  "
Set-TextValue "D51" "20.46"
$ws.Range("E51").Value = "  -4.00%  "
